$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All source cells are text (inline strings) in the original workbook.
# Force text format before assigning so Excel does not auto-convert
# numeric-looking strings (e.g. "319.43") into real numbers.
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "B35", "C35", "D35", "E35", "B36", "C36", "D36", "E36", "D37", "E37", "B38", "C38", "D38", "E38", "B39", "C39", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "B45", "C45", "D45", "E45", "B46", "C46", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cell in $cells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.935.81"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "1.905.28"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "319.43"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "0.5039"
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "0.08255"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").Value = "41.99"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").Value = "1.098"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "24.05"
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").Value = "1.916.15"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "6.387"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "7.240"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").Value = "1.010"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "91.76"
$ws.Range("E17").Value = "  -3.57%  "
$ws.Range("D18").Value = "0.00001093"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").Value = "0.06498"
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("D20").Value = "18.03"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "5.922"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "29.990.76"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "11.27"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "2.202"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "22.23"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("D27").Value = "2.135.83"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "162.04"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("D29").Value = "2.265"
$ws.Range("E29").Value = "  -5.63%  "
$ws.Range("D30").Value = "128.76"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").Value = "1.114"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "0.1036"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("D33").Value = "5.943"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("D34").Value = "3.803"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "5.375"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.02433"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("D37").Value = "0.06333"
$ws.Range("E37").Value = "  -3.88%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2147"
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "0.6577"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").Value = "1.190"
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("D41").Value = "8.667"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").Value = "11.31"
$ws.Range("E42").Value = "  -5.39%  "
$ws.Range("D43").Value = "1.206"
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("D44").Value = "2.195"
$ws.Range("E44").Value = "  +6.32%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.33"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6035"
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("D47").Value = "3.632"
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").Value = "122.79"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("D49").Value = "1.208"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").Value = "78.29"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").Value = "1.131"
$ws.Range("E51").Value = "  -2.62%  "
